$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - column F "想去人数" (number of people wanting to go) updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 5482
$ws1.Range("F3").Value = 598
$ws1.Range("F4").Value = 12024
$ws1.Range("F5").Value = 296
$ws1.Range("F6").Value = 608
$ws1.Range("F7").Value = 178
$ws1.Range("F8").Value = 315
$ws1.Range("F9").Value = 1096
$ws1.Range("F10").Value = 104

# Sheet "演出" (performances) - column F updates
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 26
$ws2.Range("F7").Value = 2

# Sheet "全部类型" (all types) - column F updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 5482
$ws4.Range("F5").Value = 598
$ws4.Range("F6").Value = 26
$ws4.Range("F7").Value = 12024
$ws4.Range("F8").Value = 296
$ws4.Range("F9").Value = 608
$ws4.Range("F10").Value = 178
$ws4.Range("F13").Value = 315
$ws4.Range("F14").Value = 1096
$ws4.Range("F15").Value = 2
$ws4.Range("F16").Value = 104
